# Updated site info and param list
#
# - The parameter "FRAC_GROWTHRESP" moves out of the "Photosynthesis" group
#   (where it sat between E_KmO/fpseudo/fpsir and fQ) and becomes the first
#   entry of the group formerly called "Post GPP".
# - That group ("HCRIT_LITTER" ... "MOIST_COEFF") is renamed from
#   "Post GPP" to "Post C uptake".
# - The view scrolls back to the top and the selection moves to C46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "FRAC_GROWTHRESP" row from the Photosynthesis block (row 17);
# everything below shifts up one row.
$ws.Rows(17).Delete()

# Insert a new row at the top of the old "Post GPP" block (now row 32) and
# put FRAC_GROWTHRESP there, under the renamed group.
$ws.Rows(32).Insert()
$ws.Range("A32").Value2 = "FRAC_GROWTHRESP"
$ws.Range("B32").Value2 = "Post C uptake"

# Rename the remaining rows of that group from "Post GPP" to "Post C uptake".
for ($r = 33; $r -le 46; $r++) {
    $ws.Range("B$r").Value2 = "Post C uptake"
}

# Update the sheet view: scroll back to the top-left and select C46.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C46").Select()
